$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Date property (row 8, column B)
$meta.Range("B8").Value = "2026-01-16T13:49:34+00:00"

# Description property (row 12, column B)
$meta.Range("B12").Value = "Dose d'antigène"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Row 4 (fr-lm-dose-antigene.doseAntigene): Short & Definition columns (L, M)
# now reuse the same text as the corrected Description ("Dose d'antigène")
$elements.Range("L4").Value = "Dose d'antigène"
$elements.Range("M4").Value = "Dose d'antigène"
